# NIT-9004583876.xlsx — "Actualiza base de datos EC y agrega parte 1 de
# nuevos estado de cuenta"
#
# The "Periodo Mora" list (rows 16-42, column E) was re-sorted from
# descending (2001 .. 1711) to ascending (1711 .. 2001) order, and the
# "Valor Mora" (column F) amount that belongs to period 2001 (22400,
# as opposed to the flat 32000 for every other period) moves together
# with its period label — from row 16 down to row 42.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periods = @(
    "1711", "1712",
    "1801", "1802", "1803", "1804", "1805", "1806",
    "1807", "1808", "1809", "1810", "1811", "1812",
    "1901", "1902", "1903", "1904", "1905", "1906",
    "1907", "1908", "1909", "1910", "1911", "1912",
    "2001"
)

$firstRow = 16
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = $firstRow + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]
    if ($periods[$i] -eq "2001") {
        $ws.Cells.Item($row, 6).Value = 22400
    } else {
        $ws.Cells.Item($row, 6).Value = 32000
    }
}
